$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.152.49'
$ws.Range('E2').Value = '  -1.54%  '
$ws.Range('D3').Value = '1.832.86'
$ws.Range('E3').Value = '  -2.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9997'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4656'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.64%  '
$ws.Range('E8').Value = '  -6.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06281'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.06%  '
$ws.Range('D10').Value = '1.821.09'
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07400'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.04'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.901'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '83.40'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6185'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.26%  '
$ws.Range('D16').Value = '30.074.41'
$ws.Range('E16').Value = '  -1.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.000'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '226.21'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007292'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9988'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').Value = '2.060.17'
$ws.Range('E22').Value = '  -3.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.855'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.91%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.869'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.157'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.48'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.860'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.82%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1018'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('E30').Value = '  -2.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.057'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.782'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.79%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04775'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.132'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7069'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.685'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01821'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.609'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8951'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.931'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.18%  '
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '103.44'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.474'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4002'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.990'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1194'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '59.67'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.452'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '32.74'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.85%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.376'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.35%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05514'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.69%  '
